$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.903.25'
$ws.Range("E2").Value = '  -0.26%  '

# Row 3
$ws.Range("D3").Value = '1.897.29'
$ws.Range("E3").Value = '  +0.00%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7923'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  -4.39%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '243.93'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +0.90%  '

# Row 7
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3165'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -3.18%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.42'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  -3.87%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07232'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +2.98%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08105'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +0.20%  '

# Row 12
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7677'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +1.04%  '

# Row 13
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.576'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +6.45%  '

# Row 14
$ws.Range("D14").Value = '1.878.49'
$ws.Range("E14").Value = '  -1.07%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.69'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +0.63%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.170'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +5.67%  '

# Row 17
$ws.Range("D17").Value = '29.896.05'
$ws.Range("E17").Value = '  -0.29%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.95'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  -0.79%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.80'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +0.61%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007801'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.79%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.220'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = '  +18.43%  '

# Row 22
$ws.Range("D22").Value = '2.163.38'
$ws.Range("E22").Value = '  +0.51%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = '  -0.08%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.001'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -0.01%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1679'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  -2.38%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.490'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +2.65%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.21'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  -0.73%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.73'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  -0.75%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.066'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = '  -1.09%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.400'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +2.89%  '

# Row 31
$ws.Range("E31").Value = '  +2.55%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.496'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +5.32%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05582'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  -5.45%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.100'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = '  +0.99%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.282'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +1.51%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7422'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +1.73%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9966'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = '  -0.33%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.631'
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = '  -3.32%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01934'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +1.12%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.783'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  +0.32%  '

# Row 41
$ws.Range("D41").Value = '1.158.49'
$ws.Range("E41").Value = '  +17.25%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.46'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +2.98%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4427'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -0.07%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.934'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  +1.51%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8532'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = '  +0.13%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '104.75'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +2.77%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.000'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  -0.07%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.886'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -0.40%  '

# Row 49
$ws.Range("E49").Value = '  +1.96%  '

# Row 50
$ws.Range("B50").Value = 'SynthetixNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.039'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  +11.37%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.468'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  -0.88%  '
